$wb = $excel.ActiveWorkbook

# zh-cn sheet: update handoff/handback datetime for the 8a3b8d57... row (row 5)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-24 09:58:04"
$wsZhCn.Range("G5").Value = "2016-02-24 09:58:47"

# de-de sheet: update handoff/handback datetime for the 8a3b8d57... row (row 5)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-24 09:58:16"
$wsDeDe.Range("G5").Value = "2016-02-24 09:59:10"
